$p = $ppt.ActivePresentation

# ------------------------------------------------------------------
# Slide 1 - "Book" details slide: refresh Id/Description/ImageUrl for
# the new book (412, "King In Black" -> description/imageUrl refresh).
# ------------------------------------------------------------------
$s1 = $p.Slides.Item(1)

# Shape 2 (TextBox3): "Id: 368" -> "Id: 412", color F45805 -> D07886
$idShape = $s1.Shapes.Item(2)
$idShape.TextFrame.TextRange.Text = "Id: 412"
$idShape.TextFrame.TextRange.Font.Color.RGB = 8812752

# Shape 4 (TextBox5): Description text refresh
$descShape = $s1.Shapes.Item(4)
$descShape.TextFrame.TextRange.Text = "Description: This is King In Black"

# Shape 5 (TextBox6): ImageUrl text refresh
$imgShape = $s1.Shapes.Item(5)
$imgShape.TextFrame.TextRange.Text = "ImageUrl: https://colorless-shrimp-958.convex.cloud/api/storage/bc963ec4-0ad1-4ef3-8ed2-20c670e2359f"

# ------------------------------------------------------------------
# Slide 2 - "Chapter: 0" slide: refresh Book Id/Title/Content for the
# first chapter of king in black.
# ------------------------------------------------------------------
$s2 = $p.Slides.Item(2)

# Shape 2 (TextBox3): "Book Id: 368" -> "Book Id: 412", color F45805 -> D07886
$bookIdShape = $s2.Shapes.Item(2)
$bookIdShape.TextFrame.TextRange.Text = "Book Id: 412"
$bookIdShape.TextFrame.TextRange.Font.Color.RGB = 8812752

# Shape 3 (TextBox4): Title text refresh
$titleShape = $s2.Shapes.Item(3)
$titleShape.TextFrame.TextRange.Text = "Title: First chapter of king in black"

# Shape 4 (TextBox5): Content text refresh
$contentShape = $s2.Shapes.Item(4)
$contentShape.TextFrame.TextRange.Text = "Content: THIS IS KING IN BLACK"

# ------------------------------------------------------------------
# New Slide 3 - "Chapter: 0" the second chapter of king in black, added
# after slide 2, using the same blank layout (slideLayout7 / ppLayoutBlank).
# ------------------------------------------------------------------
$s3 = $p.Slides.Add($p.Slides.Count + 1, 12)

$chapterBox = $s3.Shapes.AddTextbox(1, 50, 50, 600, 50)
$chapterBox.Name = "TextBox2"
$chapterBox.Line.BeginArrowheadStyle = 1
$chapterBox.Line.EndArrowheadStyle = 1
$chapterBox.TextFrame.WordWrap = 1
$chapterBox.TextFrame.AutoSize = 1
$chapterBox.TextFrame.TextRange.Text = "Chapter: 0"
$chapterBox.TextFrame.TextRange.Font.Size = 28
$chapterBox.TextFrame.TextRange.Font.Bold = 1
$chapterBox.TextFrame.TextRange.Font.Italic = 0
$chapterBox.TextFrame.TextRange.Font.Color.RGB = 0

$bookIdBox3 = $s3.Shapes.AddTextbox(1, 50, 100, 600, 50)
$bookIdBox3.Name = "TextBox3"
$bookIdBox3.Line.BeginArrowheadStyle = 1
$bookIdBox3.Line.EndArrowheadStyle = 1
$bookIdBox3.TextFrame.WordWrap = 1
$bookIdBox3.TextFrame.AutoSize = 1
$bookIdBox3.TextFrame.TextRange.Text = "Book Id: 412"
$bookIdBox3.TextFrame.TextRange.Font.Size = 18
$bookIdBox3.TextFrame.TextRange.Font.Bold = 0
$bookIdBox3.TextFrame.TextRange.Font.Italic = 0
$bookIdBox3.TextFrame.TextRange.Font.Color.RGB = 8812752

$titleBox3 = $s3.Shapes.AddTextbox(1, 50, 150, 600, 50)
$titleBox3.Name = "TextBox4"
$titleBox3.Line.BeginArrowheadStyle = 1
$titleBox3.Line.EndArrowheadStyle = 1
$titleBox3.TextFrame.WordWrap = 1
$titleBox3.TextFrame.AutoSize = 1
$titleBox3.TextFrame.TextRange.Text = "Title: Second chapter of king in black"
$titleBox3.TextFrame.TextRange.Font.Size = 24
$titleBox3.TextFrame.TextRange.Font.Bold = 1
$titleBox3.TextFrame.TextRange.Font.Italic = 0
$titleBox3.TextFrame.TextRange.Font.Color.RGB = 0

$contentBox3 = $s3.Shapes.AddTextbox(1, 50, 200, 600, 50)
$contentBox3.Name = "TextBox5"
$contentBox3.Line.BeginArrowheadStyle = 1
$contentBox3.Line.EndArrowheadStyle = 1
$contentBox3.TextFrame.WordWrap = 1
$contentBox3.TextFrame.AutoSize = 1
$contentBox3.TextFrame.TextRange.Text = "Content: THIS IS KING IN BLACK"
$contentBox3.TextFrame.TextRange.Font.Size = 18
$contentBox3.TextFrame.TextRange.Font.Bold = 0
$contentBox3.TextFrame.TextRange.Font.Italic = 0
$contentBox3.TextFrame.TextRange.Font.Color.RGB = 0
